# Add an "Outlier Removal" worksheet, modeled on the existing "Missing
# Values" sheet (same layout/styling), populated with the outlier-removal
# workflow steps, placed as the last (now active) tab.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Missing Values")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Outlier Removal"

$ws.Range("A1").Value = "Action"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Content"

$ws.Range("A2").Value = "Upload CSV"
$ws.Range("B2").Value = "5 min"
$ws.Range("C2").Value = "df = pd.read_csv('file.csv')"

$ws.Range("A3").Value = "Statistical Summary"
$ws.Range("B3").Value = "2 min"
$ws.Range("C3").Value = "df.describe() to identify outliers"

$ws.Range("A4").Value = "Visualize Outliers"
$ws.Range("B4").Value = "10 min"
$ws.Range("C4").Value = "Use sns.boxplot(data=df)"

$ws.Range("A5").Value = "Remove Outliers"
$ws.Range("B5").Value = "5 min"
$ws.Range("C5").Value = "Use IQR or Z-score method to filter out outliers"

$ws.Range("A6").Value = "Verify Changes"
$ws.Range("B6").Value = "1 min"
$ws.Range("C6").Value = "Replot with sns.boxplot(data=df)"

$ws.Range("A7").Value = "Overall"
$ws.Range("B7").Value = "23 min"

$c7 = $ws.Range("C7")
$c7.Value = ""
$c7.Font.Size = 12
$c7.Font.Color = 0

$ws.Range("A1:C7").Select()
